# Applies the scraper re-run update to the Primera RFEF Group 1 2023-2024 sheet:
#  - A handful of same-kickoff-time fixture rows had their home/away data
#    (columns F:V) shuffled between rows during the re-scrape; restore the
#    corrected row-to-fixture alignment by swapping/rotating those columns.
#  - Four newly scraped fixtures are appended as rows 112-115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2) {
    # Swap columns F:V (home..url) between two rows; A (Indice) and E (data_partida)
    # stay put since they are unaffected by the diff.
    for ($c = 6; $c -le 22; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Simple pairwise swaps (row N <-> row N+1)
Swap-RowData 2 3
Swap-RowData 12 13
Swap-RowData 14 15
Swap-RowData 24 25
Swap-RowData 26 27
Swap-RowData 30 31
Swap-RowData 34 35
Swap-RowData 36 37
Swap-RowData 72 73
Swap-RowData 92 93
Swap-RowData 96 97
Swap-RowData 106 107

# Three-way rotation: new19 = old20, new20 = old21, new21 = old19
$vals = @{}
foreach ($r in 19, 20, 21) {
    $row = @()
    for ($c = 6; $c -le 22; $c++) {
        $row += ,($ws.Cells.Item($r, $c).Value2)
    }
    $vals[$r] = $row
}
function Set-RowData($r, $row) {
    for ($c = 6; $c -le 22; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 6]
    }
}
Set-RowData 19 $vals[20]
Set-RowData 20 $vals[21]
Set-RowData 21 $vals[19]

# Append four newly scraped fixtures as rows 112-115, copying formatting from
# the last existing data row (111).
$ws.Range("A111:V111").Copy($ws.Range("A112:V115"))

function Set-Fixture($row, $indice, $date, $home, $homeGoals, $away, $awayGoals, `
    $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt, `
    $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt, `
    $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt, $url) {

    $ws.Cells.Item($row, 1).Value = $indice
    $ws.Cells.Item($row, 2).Value = "spain"
    $ws.Cells.Item($row, 3).Value = "primera-rfef-group-1"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $date
    $ws.Cells.Item($row, 6).Value = $home
    $ws.Cells.Item($row, 7).Value = $homeGoals
    $ws.Cells.Item($row, 8).Value = $away
    $ws.Cells.Item($row, 9).Value = $awayGoals
    $ws.Cells.Item($row, 10).Value = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value = $homeOpenDt
    $ws.Cells.Item($row, 12).Value = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value = $homeCloseDt
    $ws.Cells.Item($row, 14).Value = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value = $drawOpenDt
    $ws.Cells.Item($row, 16).Value = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value = $drawCloseDt
    $ws.Cells.Item($row, 18).Value = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value = $awayOpenDt
    $ws.Cells.Item($row, 20).Value = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value = $awayCloseDt
    $ws.Cells.Item($row, 22).Value = $url
}

Set-Fixture 112 111 45241.70833333334 "Cornella" 3 "R. Union" 1 `
    3.09 "09/11/2023 09:12" 2.89 "11/11/2023 16:54" `
    3.06 "09/11/2023 09:12" 3.05 "11/11/2023 16:51" `
    2.26 "09/11/2023 09:12" 2.61 "11/11/2023 16:54" `
    "https://www.betexplorer.com/football/spain/primera-rfef-group-1/cornella-real-union/ltwZ9myq/"

Set-Fixture 113 112 45241.72916666666 "Barcelona B" 4 "Sabadell" 1 `
    1.79 "09/11/2023 09:12" 1.53 "11/11/2023 16:57" `
    3.28 "09/11/2023 09:12" 4.12 "11/11/2023 16:57" `
    4.16 "09/11/2023 09:12" 6.19 "11/11/2023 16:57" `
    "https://www.betexplorer.com/football/spain/primera-rfef-group-1/barcelona-sabadell/nXj7ETSR/"

Set-Fixture 114 113 45241.75 "R. Sociedad B" 2 "Arenteiro" 2 `
    1.78 "09/11/2023 09:12" 1.82 "11/11/2023 17:59" `
    3.21 "09/11/2023 09:12" 2.97 "11/11/2023 17:59" `
    4.4 "09/11/2023 09:12" 5.85 "11/11/2023 17:59" `
    "https://www.betexplorer.com/football/spain/primera-rfef-group-1/r-sociedad-arenteiro/OIHyU8SE/"

Set-Fixture 115 114 45241.79166666666 "Dep. La Coruna" 2 "SD Logrones" 0 `
    1.43 "09/11/2023 09:13" 1.36 "11/11/2023 18:54" `
    4 "09/11/2023 09:13" 4.69 "11/11/2023 18:54" `
    6.5 "09/11/2023 09:13" 9.53 "11/11/2023 18:54" `
    "https://www.betexplorer.com/football/spain/primera-rfef-group-1/dep-la-coruna-sd-logrones/zcys8R5e/"

Write-Host "edit complete"
